# Updates the cryptos list (price/volume columns) to the latest scraped
# values, and fixes the row order for InjectiveProtocol / EthereumClassic
# (rows 30 and 31 swapped places along with their data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values formatted as plain text (e.g. 
# "43.810.79" or "6.40"); force these cells to stay text so Excel does not
# reinterpret them as numbers/dates and strip significant trailing zeros.
$dCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D20","D21","D22","D24","D25","D27","D28","D30","D31","D32","D33","D34","D38","D40","D41","D42","D43","D44","D45","D46","D48","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.810.79'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '2.318.66'
$ws.Range("E3").Value = '  +4.05%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '97.62'
$ws.Range("E5").Value = '  +5.13%  '
$ws.Range("D6").Value = '271.82'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.626'
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("D10").Value = '45.44'
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("D11").Value = '0.0948'
$ws.Range("E11").Value = '  -2.57%  '
$ws.Range("D12").Value = '8.06'
$ws.Range("E12").Value = '  -3.45%  '
$ws.Range("D13").Value = '0.106'
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").Value = '2.654.44'
$ws.Range("E14").Value = '  +3.74%  '
$ws.Range("D15").Value = '15.56'
$ws.Range("E15").Value = '  +2.95%  '
$ws.Range("D16").Value = '0.878'
$ws.Range("E16").Value = '  +8.94%  '
$ws.Range("D17").Value = '2.316.69'
$ws.Range("E17").Value = '  +3.27%  '
$ws.Range("D18").Value = '43.744.32'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("E19").Value = '  +3.93%  '
$ws.Range("D20").Value = '6.40'
$ws.Range("E20").Value = '  +5.45%  '
$ws.Range("D21").Value = '73.35'
$ws.Range("E21").Value = '  +3.85%  '
$ws.Range("D22").Value = '240.24'
$ws.Range("E22").Value = '  +2.81%  '
$ws.Range("E23").Value = '  -3.67%  '
$ws.Range("D24").Value = '9.43'
$ws.Range("E24").Value = '  +3.06%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  +1.31%  '
$ws.Range("D27").Value = '11.38'
$ws.Range("E27").Value = '  -0.39%  '
$ws.Range("D28").Value = '3.49'
$ws.Range("E28").Value = '  -0.82%  '
$ws.Range("E29").Value = '  +1.85%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '38.18'
$ws.Range("E30").Value = '  -7.04%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '22.39'
$ws.Range("E31").Value = '  +7.03%  '
$ws.Range("D32").Value = '174.39'
$ws.Range("E32").Value = '  +1.25%  '
$ws.Range("D33").Value = '0.0911'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("D34").Value = '5.48'
$ws.Range("E34").Value = '  -0.28%  '
$ws.Range("E35").Value = '  +2.54%  '
$ws.Range("E36").Value = '  +2.66%  '
$ws.Range("E37").Value = '  -3.58%  '
$ws.Range("D38").Value = '4.40'
$ws.Range("E38").Value = '  +1.96%  '
$ws.Range("E39").Value = '  -5.37%  '
$ws.Range("D40").Value = '0.244'
$ws.Range("E40").Value = '  +10.01%  '
$ws.Range("D41").Value = '2.36'
$ws.Range("E41").Value = '  +9.02%  '
$ws.Range("D42").Value = '1.38'
$ws.Range("E42").Value = '  +19.26%  '
$ws.Range("D43").Value = '12.22'
$ws.Range("E43").Value = '  -5.65%  '
$ws.Range("D44").Value = '9.19'
$ws.Range("E44").Value = '  +9.97%  '
$ws.Range("D45").Value = '62.47'
$ws.Range("E45").Value = '  -2.08%  '
$ws.Range("D46").Value = '5.36'
$ws.Range("E46").Value = '  +0.34%  '
$ws.Range("E47").Value = '  +3.11%  '
$ws.Range("D48").Value = '100.35'
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("E49").Value = '  +0.68%  '
$ws.Range("D50").Value = '0.190'
$ws.Range("E50").Value = '  +15.63%  '
$ws.Range("D51").Value = '2.541.03'
